# Apply strikethrough formatting to the Q1-Q4 exercise block (everything from
# the "Q1" heading through the final "Output the list..." bullet of Q4),
# matching the target revision's commit ("LeapYear example Test Driven
# Development") where this block of now-superseded instructions was struck
# through. The blank separator paragraphs between questions are left alone.

$d = $word.ActiveDocument

$startPara = 8
$endPara = 58
$skip = @(25, 40, 46)

for ($i = $startPara; $i -le $endPara; $i++) {
    if ($skip -contains $i) {
        continue
    }
    $p = $d.Paragraphs.Item($i)
    $p.Range.Font.StrikeThrough = 1
}

Write-Output "Applied strikethrough to paragraphs $startPara-$endPara (excluding $($skip -join ', '))"
